# Updates cryptos list values (prices, volumes, and two row swaps)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.413.37"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.01%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.643.57"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.58%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "606.61"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.10%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "151.86"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.61%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.593"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.20%  "
$ws.Range("E9").Value = "  +2.00%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.391"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +7.81%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.70"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.56%  "
$ws.Range("E12").Value = "  -0.59%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "27.90"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.37%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.117.94"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.76%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "64.258.76"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.00%  "
$ws.Range("E16").Value = "  +3.72%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.640.66"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.02%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.27"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +8.53%  "
$ws.Range("E19").Value = "  +4.03%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "353.73"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.97%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.98"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.66%  "
$ws.Range("E22").Value = "  +0.11%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.75"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.31%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "66.87"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.45%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.76"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +14.12%  "
$ws.Range("E26").Value = "  +6.19%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.39"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +8.08%  "
$ws.Range("B28").Value = "Aptos"
$ws.Range("C28").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.23"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.14%  "
$ws.Range("B29").Value = "Kaspa"
$ws.Range("C29").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.166"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.12%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "550.45"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.17%  "
$ws.Range("E31").Value = "  +0.15%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.08"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.02%  "
$ws.Range("E33").Value = "  +8.53%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.78"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.53%  "
$ws.Range("E35").Value = "  +2.29%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "167.71"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.94%  "
$ws.Range("E37").Value = "  +8.77%  "
$ws.Range("E38").Value = "  +2.20%  "
$ws.Range("E39").Value = "  -0.06%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "19.65"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.26%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "168.61"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.48%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "40.11"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.01%  "
$ws.Range("E44").Value = "  +5.28%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0589"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.76%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "21.77"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.41%  "
$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.631"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.25%  "
$ws.Range("B48").Value = "dogwifhat"
$ws.Range("C48").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.06"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +15.89%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0248"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.16%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0968"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.69%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.51"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.55%  "
